$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")
$ws.Range("C1").Value = "Template updated 12/8/22."
$ws.Range("C1").Font.Color = 255
$ws.Range("C2").Value = "Samples updated 1/8/23"
$ws.Range("C2").Font.Color = 255
